$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 7.86724192799219
$ws.Range("D2").Value = 4.612663642908414
$ws.Range("E2").Value = 10.70358358851053
$ws.Range("F2").Value = 52.19851471137271
$ws.Range("G2").Value = 3.763626178778753
$ws.Range("J2").Value = 10.03118505128965
$ws.Range("K2").Value = 20.55611068277238
$ws.Range("M2").Value = 20.4905080802847
$ws.Range("N2").Value = 23.3986823768665
$ws.Range("B3").Value = 7.800049692867613
$ws.Range("D3").Value = 4.615628859511509
$ws.Range("E3").Value = 10.72028624106469
$ws.Range("F3").Value = 52.07622248718422
$ws.Range("G3").Value = 3.76761716080528
$ws.Range("J3").Value = 10.05121168611236
$ws.Range("K3").Value = 20.35511994747449
$ws.Range("M3").Value = 20.43403655097899
$ws.Range("N3").Value = 23.44639330031144
$ws.Range("B4").Value = 7.760411822498218
$ws.Range("D4").Value = 4.617747346421424
$ws.Range("E4").Value = 10.73190675883346
$ws.Range("F4").Value = 52.01276830214352
$ws.Range("G4").Value = 3.770194212534101
$ws.Range("J4").Value = 10.06454508184396
$ws.Range("K4").Value = 20.23669732772502
$ws.Range("M4").Value = 20.40403845406376
$ws.Range("N4").Value = 23.47760922542648
$ws.Range("B5").Value = 7.744683016677794
$ws.Range("D5").Value = 4.618685722838512
$ws.Range("E5").Value = 10.7369856967952
$ws.Range("F5").Value = 51.98984484057529
$ws.Range("G5").Value = 3.771276333196929
$ws.Range("J5").Value = 10.07023957109645
$ws.Range("K5").Value = 20.18974300128392
$ws.Range("M5").Value = 20.39299713440176
$ws.Range("N5").Value = 23.49081312914561
$ws.Range("B6").Value = 7.742097346482893
$ws.Range("D6").Value = 4.618846079368857
$ws.Range("E6").Value = 10.73784980257568
$ws.Range("F6").Value = 51.98621587616562
$ws.Range("G6").Value = 3.771457951903888
$ws.Range("J6").Value = 10.07120090974243
$ws.Range("K6").Value = 20.18202642459016
$ws.Range("M6").Value = 20.39123538366217
$ws.Range("N6").Value = 23.49303481087407
$ws.Range("B7").Value = 7.760197959922473
$ws.Range("D7").Value = 4.61775969744898
$ws.Range("E7").Value = 10.73197386411431
$ws.Range("F7").Value = 52.01244725627666
$ws.Range("G7").Value = 3.770208676861066
$ws.Range("J7").Value = 10.06462082249546
$ws.Range("K7").Value = 20.23605874220377
$ws.Range("M7").Value = 20.40388474762018
$ws.Range("N7").Value = 23.47778534158714
$ws.Range("B8").Value = 7.843748562110053
$ws.Range("D8").Value = 4.61362434327576
$ws.Range("E8").Value = 10.70905950077474
$ws.Range("F8").Value = 52.15394013798126
$ws.Range("G8").Value = 3.764976073697225
$ws.Range("J8").Value = 10.03787517589671
$ws.Range("K8").Value = 20.48580494137037
$ws.Range("M8").Value = 20.47007107689674
$ws.Range("N8").Value = 23.41473438329049
$ws.Range("B9").Value = 8.019650936809594
$ws.Range("D9").Value = 4.607870459805122
$ws.Range("E9").Value = 10.67494474949492
$ws.Range("F9").Value = 52.52319158497544
$ws.Range("G9").Value = 3.755713621067833
$ws.Range("J9").Value = 9.993643040003715
$ws.Range("K9").Value = 21.01279975497967
$ws.Range("M9").Value = 20.63658269909333
$ws.Range("N9").Value = 23.30633318162118
$ws.Range("B10").Value = 8.155149153147294
$ws.Range("D10").Value = 4.605069473453852
$ws.Range("E10").Value = 10.65646241772717
$ws.Range("F10").Value = 52.84956307451917
$ws.Range("G10").Value = 3.749509430295144
$ws.Range("J10").Value = 9.966138539248067
$ws.Range("K10").Value = 21.41942752657637
$ws.Range("M10").Value = 20.78072226447864
$ws.Range("N10").Value = 23.23597884928486
$ws.Range("B11").Value = 8.217909934284817
$ws.Range("D11").Value = 4.604103015025155
$ws.Range("E11").Value = 10.64947995974926
$ws.Range("F11").Value = 53.00975550474634
$ws.Range("G11").Value = 3.746815785653875
$ws.Range("J11").Value = 9.954706871520404
$ws.Range("K11").Value = 21.60791037648869
$ws.Range("M11").Value = 20.85088149748173
$ws.Range("N11").Value = 23.20598929959545
$ws.Range("B12").Value = 8.241816592297113
$ws.Range("D12").Value = 4.603781129416048
$ws.Range("E12").Value = 10.64704048192359
$ws.Range("F12").Value = 53.07207707551071
$ws.Range("G12").Value = 3.74581414586904
$ws.Range("J12").Value = 9.950533082950736
$ws.Range("K12").Value = 21.6797260623618
$ws.Range("M12").Value = 20.87809465012644
$ws.Range("N12").Value = 23.19492286794269
$ws.Range("B13").Value = 8.23666194869176
$ws.Range("D13").Value = 4.603848494777592
$ws.Range("E13").Value = 10.64755677169653
$ws.Range("F13").Value = 53.05858161138346
$ws.Range("G13").Value = 3.746029051115884
$ws.Range("J13").Value = 9.951425086542441
$ws.Range("K13").Value = 21.66424064377993
$ws.Range("M13").Value = 20.87220535062571
$ws.Range("N13").Value = 23.19729332239782
$ws.Range("B14").Value = 8.219874052875113
$ws.Range("D14").Value = 4.60407565054701
$ws.Range("E14").Value = 10.64927516334672
$ws.Range("F14").Value = 53.01484962239407
$ws.Range("G14").Value = 3.746733012357975
$ws.Range("J14").Value = 9.954360383440347
$ws.Range("K14").Value = 21.61381021198816
$ws.Range("M14").Value = 20.85310749282847
$ws.Range("N14").Value = 23.20507304511257
$ws.Range("B15").Value = 8.209608661493096
$ws.Range("D15").Value = 4.604220527329013
$ws.Range("E15").Value = 10.65035436634373
$ws.Range("F15").Value = 52.98827791066959
$ws.Range("G15").Value = 3.747166599716867
$ws.Range("J15").Value = 9.956178534745526
$ws.Range("K15").Value = 21.58297572129473
$ws.Range("M15").Value = 20.84149309996666
$ws.Range("N15").Value = 23.20987611777209
$ws.Range("B16").Value = 8.151068387193549
$ws.Range("D16").Value = 4.605138812306804
$ws.Range("E16").Value = 10.65694739060513
$ws.Range("F16").Value = 52.83932790317307
$ws.Range("G16").Value = 3.749688045435761
$ws.Range("J16").Value = 9.966907345791116
$ws.Range("K16").Value = 21.40717491255759
$ws.Range("M16").Value = 20.77622832187934
$ws.Range("N16").Value = 23.23797930079058
$ws.Range("B17").Value = 8.115428403363932
$ws.Range("D17").Value = 4.605780851720923
$ws.Range("E17").Value = 10.66135681821611
$ws.Range("F17").Value = 52.75093728822799
$ws.Range("G17").Value = 3.751267742809866
$ws.Range("J17").Value = 9.97376565986419
$ws.Range("K17").Value = 21.30018015535631
$ws.Range("M17").Value = 20.7373559137258
$ws.Range("N17").Value = 23.25573592037756
$ws.Range("B18").Value = 8.095036087923933
$ws.Range("D18").Value = 4.606179109598127
$ws.Range("E18").Value = 10.6640271765042
$ws.Range("F18").Value = 52.70120261339166
$ws.Range("G18").Value = 3.752188461288932
$ws.Range("J18").Value = 9.977812072851048
$ws.Range("K18").Value = 21.23897337954497
$ws.Range("M18").Value = 20.71543045143808
$ws.Range("N18").Value = 23.26613870490418
$ws.Range("B19").Value = 8.088150586613118
$ws.Range("D19").Value = 4.606318934055428
$ws.Range("E19").Value = 10.66495436743264
$ws.Range("F19").Value = 52.68455388056116
$ws.Range("G19").Value = 3.752502285681582
$ws.Range("J19").Value = 9.979199590744749
$ws.Range("K19").Value = 21.21830905241732
$ws.Range("M19").Value = 20.70808164796985
$ws.Range("N19").Value = 23.26969347499793
$ws.Range("B20").Value = 8.11921143139673
$ws.Range("D20").Value = 4.605709508015861
$ws.Range("E20").Value = 10.66087354304825
$ws.Range("F20").Value = 52.76023240744509
$ws.Range("G20").Value = 3.751098327983349
$ws.Range("J20").Value = 9.97302505736773
$ws.Range("K20").Value = 21.31153583926026
$ws.Range("M20").Value = 20.74144924482459
$ws.Range("N20").Value = 23.25382606878894
$ws.Range("B21").Value = 8.224801417658401
$ws.Range("D21").Value = 4.604007734022106
$ws.Range("E21").Value = 10.64876487924325
$ws.Range("F21").Value = 53.02764992925796
$ws.Range("G21").Value = 3.746525743852525
$ws.Range("J21").Value = 9.953494006874886
$ws.Range("K21").Value = 21.62861137655139
$ws.Range("M21").Value = 20.85869960260843
$ws.Range("N21").Value = 23.20278008077157
$ws.Range("B22").Value = 8.294619908219317
$ws.Range("D22").Value = 4.603152452227944
$ws.Range("E22").Value = 10.64204372277989
$ws.Range("F22").Value = 53.21208574550962
$ws.Range("G22").Value = 3.743644403053371
$ws.Range("J22").Value = 9.941633492628201
$ws.Range("K22").Value = 21.83838173388495
$ws.Range("M22").Value = 20.93908375046257
$ws.Range("N22").Value = 23.17110879733252
$ws.Range("B23").Value = 8.257289303750705
$ws.Range("D23").Value = 4.603585475220664
$ws.Range("E23").Value = 10.64552192544898
$ws.Range("F23").Value = 53.11277389262061
$ws.Range("G23").Value = 3.745172467427202
$ws.Range("J23").Value = 9.94788100943598
$ws.Range("K23").Value = 21.72621143198055
$ws.Range("M23").Value = 20.89584278056583
$ws.Range("N23").Value = 23.18785762773189
$ws.Range("B24").Value = 8.117500818886423
$ws.Range("D24").Value = 4.605741671724945
$ws.Range("E24").Value = 10.66109161023698
$ws.Range("F24").Value = 52.7560267111029
$ws.Range("G24").Value = 3.751174881391135
$ws.Range("J24").Value = 9.973359561433561
$ws.Range("K24").Value = 21.30640097688649
$ws.Range("M24").Value = 20.73959733180489
$ws.Range("N24").Value = 23.25468890760811
$ws.Range("B25").Value = 7.970884110660904
$ws.Range("D25").Value = 4.609175865234634
$ws.Range("E25").Value = 10.68301670771081
$ws.Range("F25").Value = 52.41355195156139
$ws.Range("G25").Value = 3.758113262122814
$ws.Range("J25").Value = 10.00473106873429
$ws.Range("K25").Value = 20.86657628550121
$ws.Range("M25").Value = 20.5876611483778
$ws.Range("N25").Value = 23.33402732298513
